# Append a new data row (row 91) to the "Prices" sheet, matching the
# existing convention in the sheet where every data cell is stored as
# literal text (inlineStr) rather than a typed number/date.
#
# We assign each value with a leading apostrophe so Excel stores it as
# text even though it "looks like" a number/date, then clear the style
# back to "Normal" so no number-format / quote-prefix styling sticks to
# the cell (matching the unstyled cells used throughout the rest of the
# sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prices")

$ws.Range("A91").Value = "'2025-05-31"
$ws.Range("B91").Value = "'35.5"
$ws.Range("C91").Value = "'35.01"
$ws.Range("D91").Value = "'0.94"
$ws.Range("E91").Value = "'0.253"
$ws.Range("F91").Value = "'0.09"
$ws.Range("G91").Value = "'5,352"
$ws.Range("H91").Value = "'8,013"
$ws.Range("I91").Value = "'8,063"
$ws.Range("J91").Value = "'7.2065"

$ws.Range("A91:J91").Style = "Normal"
